$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(22, 2.136175827651134, 0.4808972786212289, -4.383165293726427),
    @(23, 2.217999173306856, 0.9689598110900888, -3.176927756989126),
    @(24, -0.9483953673264298, 0.4404845237731658, 0.1736381053924823),
    @(25, -6.509960667840401, -3.707480655319409, 5.622946441858697),
    @(26, -4.197874480965453, 6.619404304986809, 1.492906965058372),
    @(27, 0.2072237714611235, 1.827899620450818, -3.208558189457861),
    @(28, 4.348754659466381, -6.65085848446549, -1.474546417765225),
    @(29, 6.596046351838384, -1.555269931239627, -1.644609145905786),
    @(30, 1.523844735375753, 2.662038494800703, -0.68764657809816),
    @(31, -5.79685175829921, 0.4115632114739198, 5.635812402799219)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
